$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue "D2" "62.171.43"
Set-TextValue "E2" "  -0.35%  "
Set-TextValue "D3" "2.444.71"
Set-TextValue "E3" "  +0.63%  "
Set-TextValue "E4" "  -0.10%  "
Set-TextValue "D5" "582.48"
Set-TextValue "E5" "  +2.07%  "
Set-TextValue "D6" "142.94"
Set-TextValue "E6" "  -0.09%  "
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "D8" "0.531"
Set-TextValue "E8" "  +0.36%  "
Set-TextValue "D9" "2.439.51"
Set-TextValue "E9" "  +0.56%  "
Set-TextValue "E10" "  +1.74%  "
Set-TextValue "E11" "  +2.84%  "
Set-TextValue "D13" "0.343"
Set-TextValue "E13" "  -2.30%  "
Set-TextValue "D14" "26.46"
Set-TextValue "E14" "  -0.10%  "
Set-TextValue "E15" "  +1.60%  "
Set-TextValue "E16" "  +0.66%  "
Set-TextValue "D17" "62.169.94"
Set-TextValue "E17" "  -0.10%  "
Set-TextValue "D18" "2.437.09"
Set-TextValue "E18" "  +0.63%  "
Set-TextValue "D19" "10.78"
Set-TextValue "E19" "  -2.11%  "
Set-TextValue "D20" "7.18"
Set-TextValue "E20" "  +0.82%  "
Set-TextValue "D21" "326.54"
Set-TextValue "E21" "  +0.73%  "
Set-TextValue "E22" "  -0.30%  "
Set-TextValue "E23" "  +0.00%  "
Set-TextValue "E24" "  -5.12%  "
Set-TextValue "D25" "65.73"
Set-TextValue "E25" "  +0.96%  "
Set-TextValue "D26" "9.09"
Set-TextValue "E26" "  +1.15%  "
Set-TextValue "D27" "599.96"
Set-TextValue "E27" "  -3.46%  "
Set-TextValue "D28" "0.0₃0966"
Set-TextValue "E28" "  +0.68%  "
Set-TextValue "D29" "2.566.65"
Set-TextValue "E29" "  +0.65%  "
Set-TextValue "E30" "  +0.07%  "
Set-TextValue "E31" "  -1.88%  "
Set-TextValue "D32" "7.98"
Set-TextValue "E32" "  +0.04%  "
Set-TextValue "D33" "1.89"
Set-TextValue "E33" "  +1.80%  "
Set-TextValue "D34" "0.135"
Set-TextValue "E34" "  +1.11%  "
Set-TextValue "D35" "4.88"
Set-TextValue "E35" "  -2.64%  "
Set-TextValue "E36" "  +0.19%  "
Set-TextValue "E37" "  -1.07%  "
Set-TextValue "E38" "  +0.12%  "
Set-TextValue "D39" "152.78"
Set-TextValue "E39" "  +3.87%  "
Set-TextValue "D40" "18.43"
Set-TextValue "E40" "  -0.41%  "
Set-TextValue "D41" "5.28"
Set-TextValue "E41" "  +1.27%  "
Set-TextValue "D42" "43.19"
Set-TextValue "E42" "  +2.01%  "
Set-TextValue "E43" "  -0.74%  "
Set-TextValue "E44" "  +0.05%  "
Set-TextValue "D45" "2.50"
Set-TextValue "E45" "  +1.40%  "
Set-TextValue "D46" "142.11"
Set-TextValue "E46" "  -1.71%  "
Set-TextValue "D47" "3.63"
Set-TextValue "E47" "  -1.45%  "
Set-TextValue "D48" "0.0₆0268"
Set-TextValue "E48" "  +19.99%  "
Set-TextValue "D49" "0.601"
Set-TextValue "E49" "  +1.24%  "
Set-TextValue "D50" "0.0520"
Set-TextValue "D51" "19.81"
Set-TextValue "E51" "  -0.55%  "
